$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The whole "Hora inicio"/"Hora fin" column switches from the old custom
# [$-F400]h:mm:ss AM/PM format to the plain builtin h:mm format.
# ---------------------------------------------------------------------------
$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("E4").NumberFormat = "h:mm"

# Row 5: the second meeting's start/end time moved from morning to evening
$ws.Range("D5").NumberFormat = "h:mm"
$ws.Range("E5").NumberFormat = "h:mm"
$ws.Range("D5").Value = 0.84027777777777779
$ws.Range("E5").Value = 0.90277777777777779

# ---------------------------------------------------------------------------
# Row 6 (new): Tercera reunion
# ---------------------------------------------------------------------------
$ws.Range("C6").NumberFormat = "d-mmm-yy"
$ws.Range("C6").Value = 45889

$ws.Range("D6").NumberFormat = "h:mm"
$ws.Range("D6").Value = 0.79166666666666663
$ws.Range("E6").NumberFormat = "h:mm"
$ws.Range("E6").Value = 0.875

$ws.Range("G6").WrapText = $true
$ws.Range("G6").Value = "Se realizo la Tercera reunion con el objetivo de diligenciar el script de estrategias y repartir  `nlos documentos correspondientes a administracion de configuraciones, plan de riesgos y definicion de estrategias"

$ws.Range("F6").WrapText = $true
$ws.Range("F6").Value = "Documentacion Script de estrategias`n y documentos relacionados"

$ws.Range("H6:M6").Interior.Color = 5287936

# ---------------------------------------------------------------------------
# Row 7 (new): Revision de documentos
# ---------------------------------------------------------------------------
$ws.Range("C7").NumberFormat = "d-mmm-yy"
$ws.Range("C7").Value = 45894

$ws.Range("D7").NumberFormat = "h:mm"
$ws.Range("D7").Value = 0.375
$ws.Range("E7").NumberFormat = "h:mm"
$ws.Range("E7").Value = 0.45833333333333331

$ws.Range("F7").Value = "Revision documentos "

$ws.Range("G7").WrapText = $true
$ws.Range("G7").Value = "Revision grupal de los documentos que se repartieron en la tercera reunion `n(administracion de configuraciones, plan de riesgos y definicion de estrategias)"

$ws.Range("H7:M7").Interior.Color = 5287936

# ---------------------------------------------------------------------------
# Rows 8-9: the leftover Hora inicio/fin cells go back to a plain/general
# number format (they're still empty placeholders for future meetings).
# Use PasteSpecial(formats) from a plain-styled neighbour so we reuse the
# existing "no number format" style instead of minting a new custom one.
# ---------------------------------------------------------------------------
$ws.Range("C8").Copy() | Out-Null
$ws.Range("D8:E8").PasteSpecial(-4122) | Out-Null
$ws.Range("D9:E9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 14.4
$ws.Rows.Item(4).RowHeight = 28.8
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Selection moved to F13 (single cell) from the old F16:F17
# ---------------------------------------------------------------------------
$ws.Range("F13").Select() | Out-Null
